$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A7 number format to match the datetime format used for earlier rows
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 8 with the latest data point (2021-11-12)
$ws.Range("A8").Value = 44512
$ws.Range("A8").NumberFormat = "YYYY-MM-DD"
$ws.Range("B8").Value = 55473
